# Applies the scheduled market-data refresh described in the commit diff.
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H:N) for the
# affected leve rows across all eight job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 51
$ws.Range("H51").Value = 4444.875

# Row 64
$ws.Range("H64").Value = 4307.1304
$ws.Range("I64").Value = 3950
$ws.Range("J64").Value = 4497.6
$ws.Range("K64").Value = 3950
$ws.Range("L64").Value = 4497.6
$ws.Range("M64").Value = -3702
$ws.Range("N64").Value = -4993.6

# Row 67
$ws.Range("H67").Value = 4307.1304
$ws.Range("I67").Value = 3950
$ws.Range("J67").Value = 4497.6
$ws.Range("K67").Value = 3950
$ws.Range("L67").Value = 4497.6
$ws.Range("M67").Value = -3092
$ws.Range("N67").Value = -6213.6

# Row 74
$ws.Range("H74").Value = 3941.353
$ws.Range("I74").Value = 3918.4546
$ws.Range("J74").Value = 3983.3333
$ws.Range("K74").Value = 3918.4546
$ws.Range("L74").Value = 3983.3333
$ws.Range("M74").Value = -2982.4546
$ws.Range("N74").Value = -5855.3333

# Row 76
$ws.Range("H76").Value = 3540.5264
$ws.Range("I76").Value = 3355
$ws.Range("J76").Value = 3746.6667
$ws.Range("K76").Value = 3355
$ws.Range("L76").Value = 3746.6667
$ws.Range("M76").Value = -3040
$ws.Range("N76").Value = -4376.6667

# Row 77
$ws.Range("H77").Value = 3941.353
$ws.Range("I77").Value = 3918.4546
$ws.Range("J77").Value = 3983.3333
$ws.Range("K77").Value = 19592.273
$ws.Range("L77").Value = 19916.6665
$ws.Range("M77").Value = -14912.273
$ws.Range("N77").Value = -29276.6665

# Row 79
$ws.Range("H79").Value = 3540.5264
$ws.Range("I79").Value = 3355
$ws.Range("J79").Value = 3746.6667
$ws.Range("K79").Value = 3355
$ws.Range("L79").Value = 3746.6667
$ws.Range("M79").Value = -2263
$ws.Range("N79").Value = -5930.6667

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 2658.5293
$ws.Range("I45").Value = 2668.4375
$ws.Range("K45").Value = 2668.4375
$ws.Range("M45").Value = -2291.4375

# Row 63
$ws.Range("H63").Value = 3562.4
$ws.Range("I63").Value = 2600
$ws.Range("K63").Value = 2600
$ws.Range("M63").Value = -1914

# Row 66
$ws.Range("H66").Value = 3562.4
$ws.Range("I66").Value = 2600
$ws.Range("K66").Value = 13000
$ws.Range("M66").Value = -9568

# Row 74
$ws.Range("H74").Value = 13264175
$ws.Range("I74").Value = 17929568
$ws.Range("K74").Value = 17929568
$ws.Range("M74").Value = -17928694

# Row 77
$ws.Range("H77").Value = 13264175
$ws.Range("I77").Value = 17929568
$ws.Range("K77").Value = 89647840
$ws.Range("M77").Value = -89643472

# Row 122
$ws.Range("H122").Value = 1528.7142
$ws.Range("I122").Value = 1499.6
$ws.Range("J122").Value = 1601.5
$ws.Range("K122").Value = 4498.799999999999
$ws.Range("L122").Value = 4804.5
$ws.Range("M122").Value = -2048.799999999999
$ws.Range("N122").Value = -9704.5

$ws = $wb.Worksheets.Item("BSM")
# Row 92
$ws.Range("H92").Value = 15000
$ws.Range("J92").Value = 15000
$ws.Range("L92").Value = 15000
$ws.Range("N92").Value = -19992

# Row 105
$ws.Range("H105").Value = 45456220
$ws.Range("I105").Value = 55557050
$ws.Range("J105").Value = 2500
$ws.Range("K105").Value = 55557050
$ws.Range("L105").Value = 2500
$ws.Range("M105").Value = -55555303
$ws.Range("N105").Value = -5994

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1496.4445
$ws.Range("I16").Value = 1027.75
$ws.Range("J16").Value = 1871.4
$ws.Range("K16").Value = 1027.75
$ws.Range("L16").Value = 1871.4
$ws.Range("M16").Value = -740.75
$ws.Range("N16").Value = -2445.4

# Row 31
$ws.Range("H31").Value = 2728.7551
$ws.Range("I31").Value = 1062.375
$ws.Range("J31").Value = 5865.4707
$ws.Range("K31").Value = 1062.375
$ws.Range("L31").Value = 5865.4707
$ws.Range("M31").Value = -767.375
$ws.Range("N31").Value = -6455.4707

# Row 34
$ws.Range("H34").Value = 2728.7551
$ws.Range("I34").Value = 1062.375
$ws.Range("J34").Value = 5865.4707
$ws.Range("K34").Value = 1062.375
$ws.Range("L34").Value = 5865.4707
$ws.Range("M34").Value = -860.375
$ws.Range("N34").Value = -6269.4707

# Row 62
$ws.Range("H62").Value = 2998.4614
$ws.Range("I62").Value = 2886.111
$ws.Range("J62").Value = 3251.25
$ws.Range("K62").Value = 2886.111
$ws.Range("L62").Value = 3251.25
$ws.Range("M62").Value = -2262.111
$ws.Range("N62").Value = -4499.25

# Row 65
$ws.Range("H65").Value = 2998.4614
$ws.Range("I65").Value = 2886.111
$ws.Range("J65").Value = 3251.25
$ws.Range("K65").Value = 14430.555
$ws.Range("L65").Value = 16256.25
$ws.Range("M65").Value = -11310.555
$ws.Range("N65").Value = -22496.25

# Row 113
$ws.Range("H113").Value = 1496.4445
$ws.Range("I113").Value = 1027.75
$ws.Range("J113").Value = 1871.4
$ws.Range("K113").Value = 1027.75
$ws.Range("L113").Value = 1871.4
$ws.Range("M113").Value = 1142.25
$ws.Range("N113").Value = -6211.4

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 977.4706
$ws.Range("I5").Value = 732.4286
$ws.Range("K5").Value = 2197.2858
$ws.Range("M5").Value = -2085.2858

# Row 9
$ws.Range("H9").Value = 1400
$ws.Range("I9").Value = 500
$ws.Range("K9").Value = 1500
$ws.Range("M9").Value = -1276

# Row 131
$ws.Range("H131").Value = 1635.9333
$ws.Range("I131").Value = 594.4545000000001
$ws.Range("J131").Value = 4500
$ws.Range("K131").Value = 1783.3635
$ws.Range("L131").Value = 13500
$ws.Range("M131").Value = 3256.6365
$ws.Range("N131").Value = -23580

# Row 135
$ws.Range("H135").Value = 977.4706
$ws.Range("I135").Value = 732.4286
$ws.Range("K135").Value = 6591.8574
$ws.Range("M135").Value = -4056.8574

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 36255.47
$ws.Range("I70").Value = 54814.4
$ws.Range("J70").Value = 5323.9165
$ws.Range("K70").Value = 54814.4
$ws.Range("L70").Value = 5323.9165
$ws.Range("M70").Value = -54544.4
$ws.Range("N70").Value = -5863.9165

# Row 73
$ws.Range("H73").Value = 36255.47
$ws.Range("I73").Value = 54814.4
$ws.Range("J73").Value = 5323.9165
$ws.Range("K73").Value = 54814.4
$ws.Range("L73").Value = 5323.9165
$ws.Range("M73").Value = -53878.4
$ws.Range("N73").Value = -7195.9165

# Row 80
$ws.Range("H80").Value = 3419.3076
$ws.Range("I80").Value = 2774.1667
$ws.Range("J80").Value = 3972.2856
$ws.Range("K80").Value = 2774.1667
$ws.Range("L80").Value = 3972.2856
$ws.Range("M80").Value = -1776.1667
$ws.Range("N80").Value = -5968.2856

# Row 83
$ws.Range("H83").Value = 3419.3076
$ws.Range("I83").Value = 2774.1667
$ws.Range("J83").Value = 3972.2856
$ws.Range("K83").Value = 13870.8335
$ws.Range("L83").Value = 19861.428
$ws.Range("M83").Value = -8878.833500000001
$ws.Range("N83").Value = -29845.428

# Row 132
$ws.Range("H132").Value = 108642
$ws.Range("I132").Value = 113700
$ws.Range("K132").Value = 341100
$ws.Range("M132").Value = -338570

$ws = $wb.Worksheets.Item("LTW")
# Row 10
$ws.Range("H10").Value = 400
$ws.Range("I10").Value = 400
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 400
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -260
$ws.Range("N10").ClearContents()

# Row 16
$ws.Range("H16").Value = 1523.16
$ws.Range("I16").Value = 865.6667
$ws.Range("J16").Value = 4975
$ws.Range("K16").Value = 865.6667
$ws.Range("L16").Value = 4975
$ws.Range("M16").Value = -695.6667
$ws.Range("N16").Value = -5315

# Row 40
$ws.Range("H40").Value = 4413.3335
$ws.Range("I40").Value = 3920
$ws.Range("K40").Value = 3920
$ws.Range("M40").Value = -3784

# Row 55
$ws.Range("H55").Value = 104.5
$ws.Range("I55").Value = 61.375
$ws.Range("K55").Value = 61.375
$ws.Range("M55").Value = 111.625

# Row 61
$ws.Range("H61").Value = 1575.6538
$ws.Range("I61").Value = 1567.409
$ws.Range("K61").Value = 1567.409
$ws.Range("M61").Value = -1365.409

# Row 113
$ws.Range("H113").Value = 1575.6538
$ws.Range("I113").Value = 1567.409
$ws.Range("K113").Value = 1567.409
$ws.Range("M113").Value = 602.5909999999999

$ws = $wb.Worksheets.Item("WVR")
# Row 10
$ws.Range("H10").Value = 30001.334
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 30001.334
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 30001.334
$ws.Range("M10").ClearContents()
$ws.Range("N10").Value = -30339.334

# Row 81
$ws.Range("H81").Value = 2590.4285
$ws.Range("I81").Value = 1690
$ws.Range("J81").Value = 2659.6924
$ws.Range("K81").Value = 3380
$ws.Range("L81").Value = 5319.3848
$ws.Range("M81").Value = -2319
$ws.Range("N81").Value = -7441.3848

# Row 84
$ws.Range("H84").Value = 2590.4285
$ws.Range("I84").Value = 1690
$ws.Range("J84").Value = 2659.6924
$ws.Range("K84").Value = 16900
$ws.Range("L84").Value = 26596.924
$ws.Range("M84").Value = -11596
$ws.Range("N84").Value = -37204.924
